# Updates crypto price (column D) and 1-hour volume-change (column E)
# cells to match the refreshed data pulled from coinranking.com, per the
# "Updated cryptos list ... with GitHub Actions" commit.
#
# Column D prices that look like plain numbers ("304.23", "1.020", ...) are
# forced to store as text (NumberFormat "@") so Excel keeps their exact
# textual representation (e.g. trailing zeros, precise decimals) instead of
# auto-converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.449.45"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "1.629.57"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "304.23"
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3775"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3662"
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.71"
$ws.Range("E9").Value = "  -1.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08226"
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("E11").Value = "  -3.62%  "
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("E13").Value = "  -2.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.556"
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("E15").Value = "  -2.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.268"
$ws.Range("E16").Value = "  -1.79%  "
$ws.Range("D17").Value = "1.628.11"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06979"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("E20").Value = "  -2.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.457"
$ws.Range("E21").Value = "  -1.54%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.70"
$ws.Range("E23").Value = "  -0.92%  "
$ws.Range("D24").Value = "23.452.34"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.205"
$ws.Range("E25").Value = "  +3.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.478"
$ws.Range("E26").Value = "  +2.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.36"
$ws.Range("E27").Value = "  +0.34%  "
$ws.Range("E28").Value = "  -1.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.311"
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.76"
$ws.Range("E30").Value = "  -1.39%  "
$ws.Range("D31").Value = "1.812.02"
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("E32").Value = "  -3.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.797"
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.020"
$ws.Range("E34").Value = "  +5.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.77"
$ws.Range("E35").Value = "  +4.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02788"
$ws.Range("E36").Value = "  -1.14%  "
$ws.Range("E37").Value = "  -0.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08740"
$ws.Range("E38").Value = "  -1.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.07118"
$ws.Range("E39").Value = "  -3.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.026"
$ws.Range("E40").Value = "  -2.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7053"
$ws.Range("E41").Value = "  -0.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.353"
$ws.Range("E42").Value = "  -2.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.27"
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.24"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6551"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.326"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.983"
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08020"
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.201"
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "125.87"
$ws.Range("E51").Value = "  -2.79%  "
